$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "REPORTE"

# Insert a new row above row 1 (shifts everything down)
$ws.Rows.Item(1).Insert()

# Fill new header row 1 with a copy of the header labels, mirroring row 4 (now shifted)
$ws.Range("A1").Value = "FECHA"
$ws.Range("B1").Value = "CANT"
$ws.Range("C1").Value = "DESCRIPCION"
$ws.Range("D1").Value = "TALLA"
$ws.Range("E1").Value = "NOMBRE"
$ws.Range("F1").Value = "TELEFONO"
$ws.Range("G1").Value = "CONCEPTO"
$ws.Range("H1").Value = "EFECTIVO"
$ws.Range("I1").Value = "TARJETA"
$ws.Range("J1").Value = "TOTAL DE VENTA"
$ws.Range("K1").Value = "CIERRE/DIA"
$ws.Range("L1").Value = "TOTAL"

# Freeze top row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
